$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = "dlgVehicleInsuranceMain"
$ws.Range("B7").Value = "dlgVehicleInsuranceBase"

$ws.Range("A8").Value = "dlgAutomobileInsurance"
$ws.Range("B8").Value = "dlgVehicleInsuranceBase"

$ws.Range("A9").Value = "dlgTruckInsurance"
$ws.Range("B9").Value = "dlgVehicleInsuranceBase"

$ws.Range("A10").Value = "dlgMotorcycleInsurance"
$ws.Range("B10").Value = "dlgVehicleInsuranceBase"

$ws.Range("A11").Value = "dlgCamperInsurance"
$ws.Range("B11").Value = "dlgVehicleInsuranceBase"

$ws.Columns.Item(2).ColumnWidth = 23.1666666666667
$ws.Range("B7").Select()
